$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.113.09"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.143.32"
$ws.Range("E3").Value = "  +0.78%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.98"
$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.47"
$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.134.09"
$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("E9").Value = "  -0.32%  "

$ws.Range("E10").Value = "  -0.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.95"
$ws.Range("E11").Value = "  +2.98%  "

$ws.Range("E12").Value = "  -1.22%  "

$ws.Range("E13").Value = "  -1.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.41"
$ws.Range("E14").Value = "  +0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.661.79"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.925.37"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.140.32"
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.66"
$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.57"
$ws.Range("E23").Value = "  +0.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.35"
$ws.Range("E24").Value = "  +8.57%  "

$ws.Range("E25").Value = "  -1.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.44"
$ws.Range("E26").Value = "  -0.49%  "

$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.98"
$ws.Range("E28").Value = "  +11.73%  "

$ws.Range("E29").Value = "  +9.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.72"
$ws.Range("E30").Value = "  +0.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.23"
$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("E33").Value = "  +2.21%  "

$ws.Range("E34").Value = "  +1.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0846"
$ws.Range("E35").Value = "  -2.43%  "

$ws.Range("E36").Value = "  +1.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.17"
$ws.Range("E37").Value = "  +1.70%  "

$ws.Range("E38").Value = "  -2.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.20"
$ws.Range("E39").Value = "  -5.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.44"
$ws.Range("E40").Value = "  +0.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.30"
$ws.Range("E41").Value = "  +7.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "455.34"
$ws.Range("E42").Value = "  +1.57%  "

$ws.Range("E43").Value = "  +6.27%  "

$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.921.76"
$ws.Range("E45").Value = "  +1.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.33"
$ws.Range("E46").Value = "  +12.71%  "

$ws.Range("E47").Value = "  -2.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.93"
$ws.Range("E48").Value = "  +8.31%  "

$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("E50").Value = "  +2.86%  "

$ws.Range("E51").Value = "  -0.46%  "
